# Fruta / hortaliza, semanal
# Insert two new weekly data rows (Arándano (blue), Mercado Mayorista Lo Valledor
# de Santiago) before the current row 439, pushing the existing rows 439-447
# down to 441-449.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above row 439 (shifts old 439..447 -> 441..449).
$ws.Rows.Item(439).Insert()
$ws.Rows.Item(439).Insert()

# --- New row 439 -------------------------------------------------------
$ws.Cells.Item(439, 1).Value = 6
$ws.Cells.Item(439, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(439, 3).Value = "Metropolitana"
$ws.Cells.Item(439, 4).Value = 44890
$ws.Cells.Item(439, 5).Value = 13
$ws.Cells.Item(439, 6).Value = "Fruta"
$ws.Cells.Item(439, 7).Value = 100101
$ws.Cells.Item(439, 8).Value = "Berries"
$ws.Cells.Item(439, 9).Value = 100101001
$ws.Cells.Item(439, 10).Value = "Arándano (blue)"
$ws.Cells.Item(439, 11).Value = "Sin especificar"
$ws.Cells.Item(439, 12).Value = "Especial"
$ws.Cells.Item(439, 13).Value = 1500
$ws.Cells.Item(439, 14).Value = 3600
$ws.Cells.Item(439, 15).Value = 3600
$ws.Cells.Item(439, 16).Value = 3600
$ws.Cells.Item(439, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(439, 18).Value = "Región del Maule"
$ws.Cells.Item(439, 19).Value = 1800
$ws.Cells.Item(439, 20).Value = 2

# --- New row 440 -------------------------------------------------------
$ws.Cells.Item(440, 1).Value = 6
$ws.Cells.Item(440, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(440, 3).Value = "Metropolitana"
$ws.Cells.Item(440, 4).Value = 44890
$ws.Cells.Item(440, 5).Value = 13
$ws.Cells.Item(440, 6).Value = "Fruta"
$ws.Cells.Item(440, 7).Value = 100101
$ws.Cells.Item(440, 8).Value = "Berries"
$ws.Cells.Item(440, 9).Value = 100101001
$ws.Cells.Item(440, 10).Value = "Arándano (blue)"
$ws.Cells.Item(440, 11).Value = "Sin especificar"
$ws.Cells.Item(440, 12).Value = "Primera"
$ws.Cells.Item(440, 13).Value = 3500
$ws.Cells.Item(440, 14).Value = 3000
$ws.Cells.Item(440, 15).Value = 3000
$ws.Cells.Item(440, 16).Value = 3000
$ws.Cells.Item(440, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(440, 18).Value = "Región del Maule"
$ws.Cells.Item(440, 19).Value = 1500
$ws.Cells.Item(440, 20).Value = 2
